$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values
$ws.Range("C2").Value = 12.4
$ws.Range("B3").Value = 4.5999999999999996
$ws.Range("C4").Value = 1.4

# Update column widths
$ws.Columns.Item(1).ColumnWidth = 27
$ws.Columns.Item(3).ColumnWidth = 27.25
